$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informe")

$ws.Range("B2").Value = "Proyecto TSOFT"
$ws.Range("B3").Value = "TSOFT"
$ws.Range("B4").Value = "Casos de pruebas para proyecto TSOFT"
